$wb = $excel.ActiveWorkbook

# --- Add the new "Instructions" sheet as the first sheet ---
$examName = $wb.Worksheets.Item(1).Name
$instr = $wb.Worksheets.Add()
$instr.Name = "Instructions"
$examSheet = $wb.Worksheets.Item($examName)

# --- Instructions sheet content ---
$instr.Range("A1").Value = "Instructions"
$instr.Range("A2").Value = "1. Acquire an Exam Document (of your choice) with answers"
$instr.Range("A3").Value = "2. Write your answer in the ""Answer + Notes"" column. Add notes as much as you need for revision"
$instr.Range("B4").Value = "Ensure that the first letter in this area is your answer for the tool to work successfully:"
$instr.Range("B5").Value = "For Example: A but I think B could be the answer"
$instr.Range("A6").Value = "3. Go through the whole exam, only using the ""Answer + Notes"" column, matching questions to the ""Question #"""
$instr.Range("A7").Value = "4. When you're done, go through the answers and input the correct answer in the ""Paste answers here"" area."
$instr.Range("A8").Value = "5. Congratulations! You've done it. "
$instr.Range("A10").Value = "You should see an in depth analysis of the questions you were:"
$instr.Range("A11").Value = "Right"
$instr.Range("B11").Value = "Wrong"
$instr.Range("C11").Value = "??"
$instr.Range("D11").Value = "Wrong or unsure questions"
$instr.Range("C12").Value = "This count is regardless of right/wrong"
$instr.Range("D12").Value = "Unsure are any questions with notes"
$instr.Range("A14").Value = "In my opinon, the orange section is the most insightful, and I call it the ""uncertainty percentage"""
$instr.Range("A15").Value = "This might be the most useful area since over time it should measure your changes in confidence about exams!"

# --- Column widths on the Instructions sheet ---
$instr.Columns.Item(3).ColumnWidth = 37.5703125
$instr.Columns.Item(4).ColumnWidth = 36.28515625

# --- Match the colour legend used on the exam sheet ---
$instr.Range("A11").Interior.Color = 0x50B000
$instr.Range("B11").Interior.Color = 0x0000FF
$instr.Range("C11").Interior.Color = 0xE4C441
$instr.Range("D11").Interior.Color = 0xDE8B2A

# --- Tab colour + selection/active state for the Instructions sheet ---
$instr.Tab.ColorIndex = 8
$instr.Range("C19").Select() | Out-Null

# --- Fill in sample answers / notes on the exam sheet ---
$examSheet.Range("B3").Value = "A or B"
$examSheet.Range("N3").Value = "A"
$examSheet.Range("B4").Value = "B"
$examSheet.Range("N4").Value = "B"
$examSheet.Range("B5").Value = "D"
$examSheet.Range("N5").Value = "C"

$examSheet.Range("I6").Select() | Out-Null

$instr.Activate() | Out-Null

Write-Output "done"
